$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2: update rpc-reply message-id GUID
$f2 = $ws.Range("F2").Value()
$f2 = $f2.Replace("693603c7-8833-442b-8f06-a6e0f4d75ec3", "84d55363-a5d0-4921-9f37-5d569295af92")
$ws.Range("F2").Value = $f2

# H2: update edit-config response GUID, commit response GUID, and flow-id
$h2 = $ws.Range("H2").Value()
$h2 = $h2.Replace("e747da13-7a35-4c72-bdc9-bce3256ae08a", "6d842a4b-5134-49de-a195-9b1214d61bd9")
$h2 = $h2.Replace("243a9ff4-6b36-40be-bf49-2b50bac3c760", "97066db6-2ce5-48c5-82bf-1e33c8e68141")
$h2 = $h2.Replace('nc-ext:flow-id="83"', 'nc-ext:flow-id="247"')
$ws.Range("H2").Value = $h2

# I2: update rpc-reply message-id GUID
$i2 = $ws.Range("I2").Value()
$i2 = $i2.Replace("6b081b99-ffcd-473e-9ebb-124902c944b1", "08f004b2-fcad-40af-bb61-f70b20ab28b1")
$ws.Range("I2").Value = $i2
